$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is always empty in this export; drop it so
# the remaining columns (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) shift one slot to
# the left, matching the corrected xpath extraction order.
$ws.Range("E1").EntireColumn.Delete()
